$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "_2022" automation population values to "_2023"
$ws.Range("D3").Value = "AutomationTestpopulation1_2023"
$ws.Range("D7").Value = "AutomationTestpopulation2_2023"
$ws.Range("D11").Value = "AutomationTestpopulation3_2023"

$ws.Range("F3").Value = "Automation_Test_Population_1_2023"
$ws.Range("F7").Value = "Automation_Test_Population_2_2023"
$ws.Range("F11").Value = "Automation_Test_Population_3_2023"

# Update the sheet view: drop the frozen/scrolled topLeftCell and move the selection
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F11").Select()
